$wb = $excel.ActiveWorkbook

# Hunk 0: ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1177.2667
$ws.Range("J33").Value = 3082.5
$ws.Range("L33").Value = 3082.5
$ws.Range("N33").Value = -3540.5

# Hunk 1: ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3839.9678
$ws.Range("J51").Value = 4507.875
$ws.Range("L51").Value = 4507.875
$ws.Range("N51").Value = -5475.875

# Hunk 2: ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5393.25
$ws.Range("I62").Value = 5407.6665
$ws.Range("J62").Value = 5350
$ws.Range("K62").Value = 5407.6665
$ws.Range("L62").Value = 5350
$ws.Range("M62").Value = -4783.6665
$ws.Range("N62").Value = -6598

# Hunk 3: ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 5393.25
$ws.Range("I65").Value = 5407.6665
$ws.Range("J65").Value = 5350
$ws.Range("K65").Value = 27038.3325
$ws.Range("L65").Value = 26750
$ws.Range("M65").Value = -23918.3325
$ws.Range("N65").Value = -32990

# Hunk 4: ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 57895.934
$ws.Range("I98").Value = 61959.93
$ws.Range("K98").Value = 61959.93
$ws.Range("M98").Value = -60461.93

# Hunk 5: ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1553.7894
$ws.Range("I107").Value = 1276.2142
$ws.Range("K107").Value = 1276.2142
$ws.Range("M107").Value = 643.7858000000001

# Hunk 6: ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5063.125
$ws.Range("I116").Value = 5072.143
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 5072.143
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -1630.143
$ws.Range("N116").Value = -11884

# Hunk 7: ALC row 117
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 35200
$ws.Range("J117").Value = 35200
$ws.Range("L117").Value = 35200
$ws.Range("N117").Value = -44378

# Hunk 8: ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 57895.934
$ws.Range("I122").Value = 61959.93
$ws.Range("K122").Value = 185879.79
$ws.Range("M122").Value = -183429.79

# Hunk 9: ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6063783
$ws.Range("I138").Value = 2541.625
$ws.Range("J138").Value = 7095484
$ws.Range("K138").Value = 7624.875
$ws.Range("L138").Value = 21286452
$ws.Range("M138").Value = -2484.875
$ws.Range("N138").Value = -21296732

# Hunk 10: ARM row 26
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1135.1111
$ws.Range("I26").Value = 1135.1111
$ws.Range("K26").Value = 1135.1111
$ws.Range("M26").Value = -805.1111000000001

# Hunk 11: ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15159655
$ws.Range("I32").Value = 19234452
$ws.Range("J32").Value = 24691.643
$ws.Range("K32").Value = 19234452
$ws.Range("L32").Value = 24691.643
$ws.Range("M32").Value = -19234165
$ws.Range("N32").Value = -25265.643

# Hunk 12: ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 15876880
$ws.Range("I61").Value = 19233582
$ws.Range("K61").Value = 19233582
$ws.Range("M61").Value = -19233370

# Hunk 13: ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 35756176
$ws.Range("I74").Value = 37080150
$ws.Range("K74").Value = 37080150
$ws.Range("M74").Value = -37079276

# Hunk 14: ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 35756176
$ws.Range("I77").Value = 37080150
$ws.Range("K77").Value = 185400750
$ws.Range("M77").Value = -185396382

# Hunk 15: ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2530.9644
$ws.Range("I122").Value = 1720.2941
$ws.Range("J122").Value = 3783.818
$ws.Range("K122").Value = 5160.8823
$ws.Range("L122").Value = 11351.454
$ws.Range("M122").Value = -2710.8823
$ws.Range("N122").Value = -16251.454

# Hunk 16: ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 15876880
$ws.Range("I136").Value = 19233582
$ws.Range("K136").Value = 57700746
$ws.Range("M136").Value = -57698196

# Hunk 17: BSM row 2
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 58500
$ws.Range("J2").Value = 58500
$ws.Range("L2").Value = 58500
$ws.Range("N2").Value = -58726

# Hunk 18: BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5239.125
$ws.Range("I107").Value = 4485.273
$ws.Range("J107").Value = 6897.6
$ws.Range("K107").Value = 4485.273
$ws.Range("L107").Value = 6897.6
$ws.Range("M107").Value = -2565.273
$ws.Range("N107").Value = -10737.6

# Hunk 19: CRP row 6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 916.5
$ws.Range("J6").Value = 1000
$ws.Range("L6").Value = 1000
$ws.Range("N6").Value = -1226

# Hunk 20: CRP row 19
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 448.2857
$ws.Range("I19").Value = 367.6
$ws.Range("J19").Value = 650
$ws.Range("K19").Value = 367.6
$ws.Range("L19").Value = 650
$ws.Range("M19").Value = -197.6
$ws.Range("N19").Value = -990

# Hunk 21: CRP row 24
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 448.2857
$ws.Range("I24").Value = 367.6
$ws.Range("J24").Value = 650
$ws.Range("K24").Value = 367.6
$ws.Range("L24").Value = 650
$ws.Range("M24").Value = -197.6
$ws.Range("N24").Value = -990

# Hunk 22: CRP row 116
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 891891
$ws.Range("J116").Value = 891891
$ws.Range("L116").Value = 891891
$ws.Range("N116").Value = -901069

# Hunk 23: CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2338.6667
$ws.Range("I122").Value = 2238.8572
$ws.Range("K122").Value = 6716.571599999999
$ws.Range("M122").Value = -4266.571599999999

# Hunk 24: CRP row 123
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H123").Value = 42069
$ws.Range("J123").Value = 42069
$ws.Range("L123").Value = 42069
$ws.Range("N123").Value = -51869

# Hunk 25: CRP row 124
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H124").Value = 235146.33
$ws.Range("J124").Value = 235146.33
$ws.Range("L124").Value = 235146.33
$ws.Range("N124").Value = -240056.33

# Hunk 26: CUL row 56
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 35187.145
$ws.Range("I56").Value = 35187.145
$ws.Range("K56").Value = 35187.145
$ws.Range("M56").Value = -34657.145

# Hunk 27: CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 4426.357
$ws.Range("I129").Value = 3343.25
$ws.Range("J129").Value = 5238.6875
$ws.Range("K129").Value = 10029.75
$ws.Range("L129").Value = 15716.0625
$ws.Range("M129").Value = -5029.75
$ws.Range("N129").Value = -25716.0625

# Hunk 28: GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4030.3635
$ws.Range("I80").Value = 3334.4
$ws.Range("K80").Value = 3334.4
$ws.Range("M80").Value = -2336.4

# Hunk 29: GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4030.3635
$ws.Range("I83").Value = 3334.4
$ws.Range("K83").Value = 16672
$ws.Range("M83").Value = -11680

# Hunk 30: GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3609.0667
$ws.Range("J113").Value = 3283.7
$ws.Range("L113").Value = 3283.7
$ws.Range("N113").Value = -7623.7

# Hunk 31: GSM row 116
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 68000
$ws.Range("J116").Value = 68000
$ws.Range("L116").Value = 68000
$ws.Range("N116").Value = -77178

# Hunk 32: GSM row 118
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 20000
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 20000
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 20000
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -23314

# Hunk 33: GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2385.9473
$ws.Range("I122").Value = 2177.7856
$ws.Range("J122").Value = 2968.8
$ws.Range("K122").Value = 6533.3568
$ws.Range("L122").Value = 8906.400000000001
$ws.Range("M122").Value = -4083.3568
$ws.Range("N122").Value = -13806.4

# Hunk 34: GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2243
$ws.Range("I132").Value = 2055.4285
$ws.Range("K132").Value = 6166.2855
$ws.Range("M132").Value = -3636.2855

# Hunk 35: LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4291.75
$ws.Range("I40").Value = 3935.8333
$ws.Range("K40").Value = 3935.8333
$ws.Range("M40").Value = -3799.8333

# Hunk 36: LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3706.35
$ws.Range("I122").Value = 3348.647
$ws.Range("K122").Value = 10045.941
$ws.Range("M122").Value = -7595.940999999999

# Hunk 37: LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2956.2974
$ws.Range("I136").Value = 3199.303
$ws.Range("J136").Value = 951.5
$ws.Range("K136").Value = 9597.909
$ws.Range("L136").Value = 2854.5
$ws.Range("M136").Value = -7047.909
$ws.Range("N136").Value = -7954.5

# Hunk 38: WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 25001774
$ws.Range("I122").Value = 30304524
$ws.Range("J122").Value = 3099.1428
$ws.Range("K122").Value = 90913572
$ws.Range("L122").Value = 9297.428400000001
$ws.Range("M122").Value = -90911122
$ws.Range("N122").Value = -14197.4284

# Hunk 39: WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4193.926
$ws.Range("I126").Value = 5485.722
$ws.Range("J126").Value = 1610.3334
$ws.Range("K126").Value = 16457.166
$ws.Range("L126").Value = 4831.0002
$ws.Range("M126").Value = -13987.166
$ws.Range("N126").Value = -9771.0002

# Hunk 40: WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3413.589
$ws.Range("I132").Value = 3668.4834
$ws.Range("J132").Value = 2237.1538
$ws.Range("K132").Value = 11005.4502
$ws.Range("L132").Value = 6711.4614
$ws.Range("M132").Value = -8475.450199999999
$ws.Range("N132").Value = -11771.4614
